$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 124
$ws.Range("B1").Value = 199
$ws.Range("C1").Value = 124

$ws.Range("A2").Value = 124
$ws.Range("B2").Value = 173.3999999999996
$ws.Range("C2").Value = 124
